# support_crytogram.docx — fill in the cryptogram's "clue" line.
#
# The first paragraph currently holds a string of "?" placeholders (with a
# couple of proofErr-wrapped runs) standing in for the clue text:
#     ????? | ????-??????????? | ??? ???? ???????? | ??????? ?????????? ?? ??????
# It needs to become the real clue, split across several runs (one per
# word/separator, matching how the solved puzzle's other rows are already
# chunked):
#     QUEUE | PAIR-PROGRAMMING | ADVISOR METTING | MESSAGE INSTRUCTOR OR MENTOR

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$rng = $p1.Range

# Exclude the trailing paragraph mark, then drop in the full replacement
# text in one shot (this also removes the old runs + the gramStart/gramEnd
# proofErr markers that bracketed the old "| ???" run).
$pStart = $rng.Start
$rng.End = $rng.End - 1
$rng.Text = "QUEUE | PAIR-PROGRAMMING | ADVISOR METTING | MESSAGE INSTRUCTOR OR MENTOR"

# Re-split that single run into one run per word/separator (so the XML
# chunking matches the target) by nudging formatting on each sub-range
# and immediately reverting it — the engine keeps the run split even once
# the formatting is back to identical, but a mere .Text assignment would
# otherwise get coalesced into one big run.
$segments = @(
  "QUEUE",
  " | ",
  "PAIR",
  "-",
  "PROGRAMMING",
  " | ",
  "ADVISOR",
  " ",
  "METTING",
  " | ",
  "MESSAGE",
  " ",
  "INSTRUCTOR",
  " ",
  "OR",
  " ",
  "MENTOR"
)

$pos = 0
foreach ($seg in $segments) {
  $segStart = $pStart + $pos
  $segEnd = $segStart + $seg.Length
  $pos = $pos + $seg.Length

  $sub = $d.Range($segStart, $segEnd)
  $sub.Bold = 1
  $sub.Bold = 0
}

Write-Output "Paragraph 1 now reads:" $d.Paragraphs(1).Range.Text
